$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 6.818749999999999
$ws.Range("C2").Value = 3.800899999999999
$ws.Range("D2").Value = 28.62214999999999
$ws.Range("F2").Value = 26.82205
$ws.Range("G2").Value = 26.82205
$ws.Range("K2").Value = 34.8128
$ws.Range("L2").Value = 26.822
$ws.Range("M2").Value = 7.9908
$ws.Range("N2").Value = 7.9908
$ws.Range("B3").Value = 12.033
$ws.Range("C3").Value = 4.49
$ws.Range("D3").Value = 51.577
$ws.Range("F3").Value = 30.714
$ws.Range("G3").Value = 28.98
$ws.Range("H3").Value = 1.733
$ws.Range("I3").Value = 1.733
$ws.Range("K3").Value = 32.2246
$ws.Range("L3").Value = 28.98
$ws.Range("M3").Value = 3.2446
$ws.Range("N3").Value = 3.2446
$ws.Range("B4").Value = 33.148
$ws.Range("C4").Value = 12.313
$ws.Range("D4").Value = 33.21
$ws.Range("F4").Value = 36.399
$ws.Range("G4").Value = 36.39866575342466
$ws.Range("K4").Value = 35.663
$ws.Range("L4").Value = 35.663

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 6.144800000000001
$ws.Range("C2").Value = 3.6393
$ws.Range("D2").Value = 28.50754999999999
$ws.Range("F2").Value = 26.59835
$ws.Range("G2").Value = 26.59835
$ws.Range("K2").Value = 187.4636
$ws.Range("L2").Value = 26.598
$ws.Range("M2").Value = 160.8656
$ws.Range("N2").Value = 9.087400000000001
$ws.Range("O2").Value = 151.7782
$ws.Range("B3").Value = 9.163
$ws.Range("C3").Value = 3.818
$ws.Range("D3").Value = 53.394
$ws.Range("F3").Value = 29.562
$ws.Range("G3").Value = 28.076
$ws.Range("H3").Value = 1.486
$ws.Range("I3").Value = 1.486
$ws.Range("K3").Value = 109.719
$ws.Range("L3").Value = 28.076
$ws.Range("M3").Value = 81.643
$ws.Range("N3").Value = 7.8102
$ws.Range("O3").Value = 73.83279999999999
$ws.Range("B4").Value = 19.757
$ws.Range("C4").Value = 7.505
$ws.Range("D4").Value = 44.213
$ws.Range("F4").Value = 31.867
$ws.Range("G4").Value = 31.86704657534246
$ws.Range("K4").Value = 34.873
$ws.Range("L4").Value = 31.419
$ws.Range("M4").Value = 3.454
$ws.Range("N4").Value = 3.454

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 7.648300000000001
$ws.Range("C2").Value = 4.11755
$ws.Range("D2").Value = 25.8179
$ws.Range("F2").Value = 27.0422
$ws.Range("G2").Value = 27.0422
$ws.Range("K2").Value = 33.5992
$ws.Range("L2").Value = 27.042
$ws.Range("M2").Value = 6.5572
$ws.Range("N2").Value = 6.5572
$ws.Range("B3").Value = 12.033
$ws.Range("C3").Value = 4.49
$ws.Range("D3").Value = 51.577
$ws.Range("F3").Value = 31.525
$ws.Range("G3").Value = 28.98
$ws.Range("H3").Value = 2.545
$ws.Range("I3").Value = 2.545
$ws.Range("K3").Value = 28.98
$ws.Range("L3").Value = 28.98
$ws.Range("B4").Value = 33.148
$ws.Range("C4").Value = 12.313
$ws.Range("D4").Value = 33.21
$ws.Range("F4").Value = 36.399
$ws.Range("G4").Value = 36.39866575342466
$ws.Range("K4").Value = 35.663
$ws.Range("L4").Value = 35.663

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 8.215900000000001
$ws.Range("C2").Value = 4.038349999999999
$ws.Range("D2").Value = 26.7446
$ws.Range("F2").Value = 27.2377
$ws.Range("G2").Value = 27.2377
$ws.Range("K2").Value = 29.5322
$ws.Range("L2").Value = 27.238
$ws.Range("M2").Value = 2.2942
$ws.Range("N2").Value = 2.2942
$ws.Range("B3").Value = 12.033
$ws.Range("C3").Value = 4.49
$ws.Range("D3").Value = 51.577
$ws.Range("F3").Value = 31.525
$ws.Range("G3").Value = 28.98
$ws.Range("H3").Value = 2.545
$ws.Range("I3").Value = 2.545
$ws.Range("K3").Value = 28.98
$ws.Range("L3").Value = 28.98
$ws.Range("B4").Value = 33.148
$ws.Range("C4").Value = 12.313
$ws.Range("D4").Value = 33.21
$ws.Range("F4").Value = 36.399
$ws.Range("G4").Value = 36.39866575342466
$ws.Range("K4").Value = 35.663
$ws.Range("L4").Value = 35.663

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 7.955849999999998
$ws.Range("C2").Value = 3.772699999999999
$ws.Range("D2").Value = 27.4432
$ws.Range("F2").Value = 27.1573
$ws.Range("G2").Value = 27.1573
$ws.Range("K2").Value = 30.6968
$ws.Range("L2").Value = 27.157
$ws.Range("M2").Value = 3.539800000000001
$ws.Range("N2").Value = 3.539800000000001
$ws.Range("O2").Value = 0
$ws.Range("B3").Value = 12.033
$ws.Range("C3").Value = 4.49
$ws.Range("D3").Value = 48.758
$ws.Range("F3").Value = 31.467
$ws.Range("G3").Value = 28.922
$ws.Range("H3").Value = 2.545
$ws.Range("I3").Value = 2.545
$ws.Range("K3").Value = 29.1636
$ws.Range("L3").Value = 28.922
$ws.Range("M3").Value = 0.2416
$ws.Range("N3").Value = 0.2416
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 33.148
$ws.Range("C4").Value = 12.313
$ws.Range("D4").Value = 32.702
$ws.Range("F4").Value = 36.388
$ws.Range("G4").Value = 36.3882191780822
$ws.Range("K4").Value = 35.6966
$ws.Range("L4").Value = 35.653
$ws.Range("M4").Value = 0.0436
$ws.Range("N4").Value = 0.0436
